$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 692.025
$ws.Range("I2").Value = 575.4857
$ws.Range("J2").Value = 1507.8
$ws.Range("K2").Value = 575.4857
$ws.Range("L2").Value = 1507.8
$ws.Range("M2").Value = -462.4857
$ws.Range("N2").Value = -1733.8

$ws.Range("H16").Value = 115708
$ws.Range("I16").Value = 115708
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 115708
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -115421

$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("N45").ClearContents()

$ws.Range("H61").Value = 3033814
$ws.Range("I61").Value = 3349.5557
$ws.Range("J61").Value = 6670371
$ws.Range("K61").Value = 3349.5557
$ws.Range("L61").Value = 6670371
$ws.Range("M61").Value = -3137.5557
$ws.Range("N61").Value = -6670795

$ws.Range("H80").Value = 29090
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 29090
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 29090
$ws.Range("N80").Value = -31086

$ws.Range("H83").Value = 29090
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 29090
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 87270
$ws.Range("N83").Value = -97254

$ws.Range("H116").Value = 692.025
$ws.Range("I116").Value = 575.4857
$ws.Range("J116").Value = 1507.8
$ws.Range("K116").Value = 575.4857
$ws.Range("L116").Value = 1507.8
$ws.Range("M116").Value = 1718.5143
$ws.Range("N116").Value = -6095.8

$ws.Range("H132").Value = 4074.3594
$ws.Range("I132").Value = 2738.4595
$ws.Range("J132").Value = 5905.037
$ws.Range("K132").Value = 8215.378499999999
$ws.Range("L132").Value = 17715.111
$ws.Range("M132").Value = -5685.378499999999
$ws.Range("N132").Value = -22775.111

$ws.Range("H136").Value = 3033814
$ws.Range("I136").Value = 3349.5557
$ws.Range("J136").Value = 6670371
$ws.Range("K136").Value = 10048.6671
$ws.Range("L136").Value = 20011113
$ws.Range("M136").Value = -7498.667099999999
$ws.Range("N136").Value = -20016213


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 692.025
$ws.Range("I3").Value = 575.4857
$ws.Range("J3").Value = 1507.8
$ws.Range("K3").Value = 575.4857
$ws.Range("L3").Value = 1507.8
$ws.Range("M3").Value = -461.4857
$ws.Range("N3").Value = -1735.8

$ws.Range("H20").Value = 9201.788
$ws.Range("I20").Value = 8364.308000000001
$ws.Range("J20").Value = 12312.429
$ws.Range("K20").Value = 8364.308000000001
$ws.Range("L20").Value = 12312.429
$ws.Range("M20").Value = -8117.308000000001
$ws.Range("N20").Value = -12806.429


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("M33").ClearContents()

$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

$ws.Range("H74").Value = 42761
$ws.Range("I74").Value = 42285
$ws.Range("J74").Value = 42999
$ws.Range("K74").Value = 42285
$ws.Range("L74").Value = 42999
$ws.Range("M74").Value = -41411
$ws.Range("N74").Value = -44747

$ws.Range("H77").Value = 42761
$ws.Range("I77").Value = 42285
$ws.Range("J77").Value = 42999
$ws.Range("K77").Value = 126855
$ws.Range("L77").Value = 128997
$ws.Range("M77").Value = -122487
$ws.Range("N77").Value = -137733

$ws.Range("H86").Value = 36075.758
$ws.Range("I86").Value = 163287.6
$ws.Range("J86").Value = 9573.291999999999
$ws.Range("K86").Value = 163287.6
$ws.Range("L86").Value = 9573.291999999999
$ws.Range("M86").Value = -162164.6
$ws.Range("N86").Value = -11819.292

$ws.Range("H89").Value = 36075.758
$ws.Range("I89").Value = 163287.6
$ws.Range("J89").Value = 9573.291999999999
$ws.Range("K89").Value = 816438
$ws.Range("L89").Value = 47866.46
$ws.Range("M89").Value = -810822
$ws.Range("N89").Value = -59098.46

$ws.Range("H99").Value = 1014798.5
$ws.Range("I99").Value = 24179.8
$ws.Range("K99").Value = 24179.8
$ws.Range("M99").Value = -22681.8

$ws.Range("H126").Value = 1014798.5
$ws.Range("I126").Value = 24179.8
$ws.Range("K126").Value = 72539.39999999999
$ws.Range("M126").Value = -70069.39999999999


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3408.6667
$ws.Range("I7").Value = 3000.6
$ws.Range("J7").Value = 3918.75
$ws.Range("K7").Value = 3000.6
$ws.Range("L7").Value = 3918.75
$ws.Range("M7").Value = -2888.6
$ws.Range("N7").Value = -4142.75

$ws.Range("H16").Value = 1691.6923
$ws.Range("I16").Value = 1417.3636
$ws.Range("J16").Value = 3200.5
$ws.Range("K16").Value = 1417.3636
$ws.Range("L16").Value = 3200.5
$ws.Range("M16").Value = -1247.3636
$ws.Range("N16").Value = -3540.5

$ws.Range("H22").Value = 1394
$ws.Range("I22").Value = 915
$ws.Range("J22").Value = 1633.5
$ws.Range("K22").Value = 915
$ws.Range("L22").Value = 1633.5
$ws.Range("M22").Value = -620
$ws.Range("N22").Value = -2223.5

$ws.Range("H27").Value = 1394
$ws.Range("I27").Value = 915
$ws.Range("J27").Value = 1633.5
$ws.Range("K27").Value = 915
$ws.Range("L27").Value = 1633.5
$ws.Range("M27").Value = -808
$ws.Range("N27").Value = -1847.5

$ws.Range("H46").Value = 4166.6665
$ws.Range("I46").Value = 2500
$ws.Range("J46").Value = 5000
$ws.Range("K46").Value = 2500
$ws.Range("L46").Value = 5000
$ws.Range("M46").Value = -2312
$ws.Range("N46").Value = -5376

$ws.Range("H126").Value = 3408.6667
$ws.Range("I126").Value = 3000.6
$ws.Range("J126").Value = 3918.75
$ws.Range("K126").Value = 9001.799999999999
$ws.Range("L126").Value = 11756.25
$ws.Range("M126").Value = -6531.799999999999
$ws.Range("N126").Value = -16696.25

$ws.Range("H132").Value = 4633834.5
$ws.Range("I132").Value = 5956444
$ws.Range("J132").Value = 4701
$ws.Range("K132").Value = 17869332
$ws.Range("L132").Value = 14103
$ws.Range("M132").Value = -17866802
$ws.Range("N132").Value = -19163

